$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.045859872611465
$ws.Range("C2").Value = 0.0522292993630573
$ws.Range("D2").Value = 0.973248407643312
$ws.Range("E2").Value = 0.00382165605095541
$ws.Range("F2").Value = 0.0178343949044586
$ws.Range("G2").Value = 0.0140127388535032
$ws.Range("H2").Value = 0.00254777070063694
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0.974522292993631
$ws.Range("K2").Value = 0.00254777070063694
$ws.Range("L2").Value = 0.0089171974522293
$ws.Range("M2").Value = 0.0292993630573248
$ws.Range("N2").Value = 0.629299363057325
$ws.Range("O2").Value = 0.0114649681528662
$ws.Range("P2").Value = 0.021656050955414
$ws.Range("Q2").Value = 0.971974522292994
$ws.Range("R2").Value = 0.00636942675159236
$ws.Range("S2").Value = 0.964331210191083
$ws.Range("T2").Value = 0.00254777070063694
$ws.Range("U2").Value = 0.793630573248408
$ws.Range("V2").Value = 0.0101910828025478
$ws.Range("W2").Value = 0.00127388535031847
$ws.Range("X2").Value = 0.0089171974522293

# Row 3
$ws.Range("B3").Value = 0.921019108280255
$ws.Range("C3").Value = 0.932484076433121
$ws.Range("D3").Value = 0.00509554140127389
$ws.Range("E3").Value = 0.00254777070063694
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.997452229299363
$ws.Range("I3").Value = 0.998726114649682
$ws.Range("J3").Value = 0.00127388535031847
$ws.Range("K3").Value = 0.00127388535031847
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0.0636942675159236
$ws.Range("O3").Value = 0.361783439490446
$ws.Range("P3").Value = 0.896815286624204
$ws.Range("Q3").Value = 0.00509554140127389
$ws.Range("R3").Value = 0.0114649681528662
$ws.Range("S3").Value = 0.00127388535031847
$ws.Range("T3").Value = 0.0140127388535032
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0.185987261146497
$ws.Range("W3").Value = 0.974522292993631
$ws.Range("X3").Value = 0.987261146496815

# Row 4
$ws.Range("B4").Value = 0.0127388535031847
$ws.Range("C4").Value = 0.00509554140127389
$ws.Range("D4").Value = 0.00636942675159236
$ws.Range("E4").Value = 0.992356687898089
$ws.Range("F4").Value = 0.982165605095541
$ws.Range("G4").Value = 0.984713375796178
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0.00509554140127389
$ws.Range("K4").Value = 0.993630573248408
$ws.Range("L4").Value = 0.989808917197452
$ws.Range("M4").Value = 0.969426751592357
$ws.Range("N4").Value = 0.304458598726115
$ws.Range("O4").Value = 0.0114649681528662
$ws.Range("P4").Value = 0.0089171974522293
$ws.Range("Q4").Value = 0.00764331210191083
$ws.Range("R4").Value = 0.975796178343949
$ws.Range("S4").Value = 0.0229299363057325
$ws.Range("T4").Value = 0.00764331210191083
$ws.Range("U4").Value = 0.0165605095541401
$ws.Range("V4").Value = 0.021656050955414
$ws.Range("W4").Value = 0.0229299363057325
$ws.Range("X4").Value = 0.00127388535031847

# Row 5
$ws.Range("B5").Value = 0.0191082802547771
$ws.Range("C5").Value = 0.0101910828025478
$ws.Range("D5").Value = 0.0152866242038217
$ws.Range("E5").Value = 0.00127388535031847
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.00127388535031847
$ws.Range("J5").Value = 0.0191082802547771
$ws.Range("K5").Value = 0.00254777070063694
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0.00254777070063694
$ws.Range("O5").Value = 0.615286624203822
$ws.Range("P5").Value = 0.0700636942675159
$ws.Range("Q5").Value = 0.0140127388535032
$ws.Range("R5").Value = 0.00636942675159236
$ws.Range("S5").Value = 0.0114649681528662
$ws.Range("T5").Value = 0.975796178343949
$ws.Range("U5").Value = 0.189808917197452
$ws.Range("V5").Value = 0.782165605095541
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0.00254777070063694

